$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pass 0: extend the date-formatted column A down to row 31 ---------
# Rows 13-17 already carry the "A-column date" style from the template;
# rows 18-31 are brand new, so copy that style down before writing values
# (copying formats only, so no new style / number-format entries spawn).
$ws.Range("A17").Copy()
$ws.Range("A18:A31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Pass 1: "priming" writes -----------------------------------------
# The engine appends a brand-new shared string to the sharedStrings table
# the first time a cell is set to that (not-yet-seen) text, in the exact
# order the assignments happen. To reproduce the target workbook's shared
# string table layout exactly, we first touch one cell per brand-new
# string, in the precise order the strings must appear in the table.
# (Every one of these cells also happens to be its own final value, so
# nothing here is wasted / needs to be undone.)
$ws.Cells.Item(13, 4).Value2 = "Aufgabenverteilung, Organisierung "
$ws.Cells.Item(18, 3).Value2 = "Offline-Treffen"
$ws.Cells.Item(12, 3).Value2 = "Aufgaben"
$ws.Cells.Item(23, 4).Value2 = "Importer anpassungen"
$ws.Cells.Item(18, 4).Value2 = "Importer überarbeiten"
$ws.Cells.Item(31, 4).Value2 = "Licht sphären farbe, Erfahrungsbericht, Präsentation vorbereitung"
$ws.Cells.Item(30, 4).Value2 = "Licht sphären anzeigen"
$ws.Cells.Item(29, 4).Value2 = "Bugfixing und Farbige objekte anzeigen"
$ws.Cells.Item(28, 4).Value2 = "Bugfixing von preview"
$ws.Cells.Item(21, 4).Value2 = "Indices hinzugefügt"
$ws.Cells.Item(22, 4).Value2 = "vertice / indice anpassen"
$ws.Cells.Item(20, 4).Value2 = "bugfixing"
$ws.Cells.Item(24, 4).Value2 = "Implementierung der Preview"
$ws.Cells.Item(25, 4).Value2 = "Converter hinzugefügt"
$ws.Cells.Item(26, 4).Value2 = "Converter angepasst & bugfixing"

# --- Pass 2: fill in the full A:D grid for rows 12-31 ------------------
$rows = @(
  @{ R = 12; A = 45592; B = 4;    C = "Aufgaben";         D = "Erstes konzept der Raytracer Math library" },
  @{ R = 13; A = 45593; B = 2;    C = "Besprechung";      D = "Aufgabenverteilung, Organisierung " },
  @{ R = 14; A = 45598; B = 1.5;  C = "Besprechung";      D = "Aufgabenverteilung, Organisierung " },
  @{ R = 15; A = 45600; B = 2.15; C = "Besprechung";      D = "Aufgabenverteilung, Organisierung " },
  @{ R = 16; A = 45605; B = 3;    C = "Besprechung";      D = "Aufgabenverteilung, Organisierung " },
  @{ R = 17; A = 45612; B = 1.5;  C = "Besprechung";      D = "Aufgabenverteilung, Organisierung " },
  @{ R = 18; A = 45616; B = 5;    C = "Offline-Treffen";  D = "Importer überarbeiten" },
  @{ R = 19; A = 45619; B = 1.5;  C = "Besprechung";      D = "Aufgabenverteilung, Organisierung " },
  @{ R = 20; A = 45624; B = 3;    C = "Offline-Treffen";  D = "bugfixing" },
  @{ R = 21; A = 45628; B = 6;    C = "Offline-Treffen";  D = "Indices hinzugefügt" },
  @{ R = 22; A = 45638; B = 6;    C = "Offline-Treffen";  D = "vertice / indice anpassen" },
  @{ R = 23; A = 45671; B = 6;    C = "Aufgaben";         D = "Importer anpassungen" },
  @{ R = 24; A = 45672; B = 7;    C = "Aufgaben";         D = "Implementierung der Preview" },
  @{ R = 25; A = 45673; B = 6;    C = "Aufgaben";         D = "Converter hinzugefügt" },
  @{ R = 26; A = 45674; B = 6;    C = "Aufgaben";         D = "Converter angepasst & bugfixing" },
  @{ R = 27; A = 45675; B = 4;    C = "Aufgaben";         D = "Implementierung der Preview" },
  @{ R = 28; A = 45676; B = 6;    C = "Aufgaben";         D = "Bugfixing von preview" },
  @{ R = 29; A = 45677; B = 5;    C = "Aufgaben";         D = "Bugfixing und Farbige objekte anzeigen" },
  @{ R = 30; A = 45678; B = 6;    C = "Aufgaben";         D = "Licht sphären anzeigen" },
  @{ R = 31; A = 45679; B = 4.5;  C = "Aufgaben";         D = "Licht sphären farbe, Erfahrungsbericht, Präsentation vorbereitung" }
)

foreach ($row in $rows) {
  $ws.Cells.Item($row.R, 1).Value2 = $row.A
  $ws.Cells.Item($row.R, 2).Value2 = $row.B
  $ws.Cells.Item($row.R, 3).Value2 = $row.C
  $ws.Cells.Item($row.R, 4).Value2 = $row.D
}

# --- View state: used range grew to A1:D31, scrolled & zoomed differently
$ws.Range("A13").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130
$ws.Range("B31").Select() | Out-Null
